$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 4272.125
$ws.Range("I106").Value = 4435.4
$ws.Range("J106").Value = 4000
$ws.Range("K106").Value = 4435.4
$ws.Range("L106").Value = 4000
$ws.Range("M106").Value = -3804.4
$ws.Range("N106").Value = -5262

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1463.9231
$ws.Range("I97").Value = 918.3333
$ws.Range("K97").Value = 918.3333
$ws.Range("M97").Value = -422.3333

$ws.Range("H122").Value = 1931.4375
$ws.Range("I122").Value = 1799.1538
$ws.Range("J122").Value = 2504.6667
$ws.Range("K122").Value = 5397.4614
$ws.Range("L122").Value = 7514.000100000001
$ws.Range("M122").Value = -2947.4614
$ws.Range("N122").Value = -12414.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2783
$ws.Range("I99").Value = 2766
$ws.Range("J99").Value = 2800
$ws.Range("K99").Value = 2766
$ws.Range("L99").Value = 2800
$ws.Range("M99").Value = -1268
$ws.Range("N99").Value = -5796

$ws.Range("H107").Value = 879.1739
$ws.Range("I107").Value = 866.55
$ws.Range("J107").Value = 963.3333
$ws.Range("K107").Value = 866.55
$ws.Range("L107").Value = 963.3333
$ws.Range("M107").Value = 1053.45
$ws.Range("N107").Value = -4803.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11367076
$ws.Range("I31").Value = 25001216
$ws.Range("J31").Value = 5293
$ws.Range("K31").Value = 25001216
$ws.Range("L31").Value = 5293
$ws.Range("M31").Value = -25000921
$ws.Range("N31").Value = -5883

$ws.Range("H34").Value = 11367076
$ws.Range("I34").Value = 25001216
$ws.Range("J34").Value = 5293
$ws.Range("K34").Value = 25001216
$ws.Range("L34").Value = 5293
$ws.Range("M34").Value = -25001014
$ws.Range("N34").Value = -5697

$ws.Range("H58").Value = 2547.7222
$ws.Range("I58").Value = 2101.0908
$ws.Range("J58").Value = 2744.24
$ws.Range("K58").Value = 2101.0908
$ws.Range("L58").Value = 2744.24
$ws.Range("M58").Value = -1898.0908
$ws.Range("N58").Value = -3150.24

$ws.Range("H107").Value = 1289.625
$ws.Range("I107").Value = 790.0909
$ws.Range("J107").Value = 1712.3077
$ws.Range("K107").Value = 790.0909
$ws.Range("L107").Value = 1712.3077
$ws.Range("M107").Value = 1129.9091
$ws.Range("N107").Value = -5552.3077

$ws.Range("H136").Value = 2547.7222
$ws.Range("I136").Value = 2101.0908
$ws.Range("J136").Value = 2744.24
$ws.Range("K136").Value = 6303.2724
$ws.Range("L136").Value = 8232.719999999999
$ws.Range("M136").Value = -3753.2724
$ws.Range("N136").Value = -13332.72

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3380.1333
$ws.Range("I70").Value = 1528.8572
$ws.Range("K70").Value = 4586.571599999999
$ws.Range("M70").Value = -4271.571599999999

$ws.Range("H73").Value = 3380.1333
$ws.Range("I73").Value = 1528.8572
$ws.Range("K73").Value = 4586.571599999999
$ws.Range("M73").Value = -3494.571599999999

$ws.Range("H75").Value = 659.6
$ws.Range("I75").Value = 300
$ws.Range("J75").Value = 749.5
$ws.Range("K75").Value = 900
$ws.Range("L75").Value = 2248.5
$ws.Range("M75").Value = 98
$ws.Range("N75").Value = -4244.5

$ws.Range("H78").Value = 659.6
$ws.Range("I78").Value = 300
$ws.Range("J78").Value = 749.5
$ws.Range("K78").Value = 2700
$ws.Range("L78").Value = 6745.5
$ws.Range("M78").Value = 2292
$ws.Range("N78").Value = -16729.5

$ws.Range("H129").Value = 1618.3125
$ws.Range("I129").Value = 480
$ws.Range("J129").Value = 1780.9286
$ws.Range("K129").Value = 1440
$ws.Range("L129").Value = 5342.7858
$ws.Range("M129").Value = 3560
$ws.Range("N129").Value = -15342.7858

$ws.Range("H131").Value = 787.03845
$ws.Range("J131").Value = 1009.44446
$ws.Range("L131").Value = 3028.33338
$ws.Range("N131").Value = -13108.33338

$ws.Range("H138").Value = 27778970
$ws.Range("I138").Value = 31250990
$ws.Range("J138").Value = 2800
$ws.Range("K138").Value = 93752970
$ws.Range("L138").Value = 8400
$ws.Range("M138").Value = -93747830
$ws.Range("N138").Value = -18680

$ws.Range("H139").Value = 2874.2856
$ws.Range("I139").Value = 1693.7142
$ws.Range("J139").Value = 4645.143
$ws.Range("K139").Value = 5081.142599999999
$ws.Range("L139").Value = 13935.429
$ws.Range("M139").Value = 58.85740000000078
$ws.Range("N139").Value = -24215.429

$ws.Range("H141").Value = 6771.857
$ws.Range("I141").Value = 4216.364
$ws.Range("J141").Value = 8425.412
$ws.Range("K141").Value = 12649.092
$ws.Range("L141").Value = 25276.236
$ws.Range("M141").Value = -7469.091999999999
$ws.Range("N141").Value = -35636.236

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 17323.125
$ws.Range("I9").Value = 292.5
$ws.Range("J9").Value = 23000
$ws.Range("K9").Value = 292.5
$ws.Range("L9").Value = 23000
$ws.Range("M9").Value = -68.5
$ws.Range("N9").Value = -23448

$ws.Range("H13").Value = 504.5
$ws.Range("I13").Value = 504.5
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 504.5
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -364.5
$ws.Range("N13").ClearContents()

$ws.Range("H16").Value = 1012.9048
$ws.Range("I16").Value = 774.64703
$ws.Range("J16").Value = 2025.5
$ws.Range("K16").Value = 774.64703
$ws.Range("L16").Value = 2025.5
$ws.Range("M16").Value = -604.64703
$ws.Range("N16").Value = -2365.5

$ws.Range("H132").Value = 2262.9412
$ws.Range("I132").Value = 1640.7858
$ws.Range("J132").Value = 5166.3335
$ws.Range("K132").Value = 4922.357400000001
$ws.Range("L132").Value = 15499.0005
$ws.Range("M132").Value = -2392.357400000001
$ws.Range("N132").Value = -20559.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 562.8095
$ws.Range("I107").Value = 482.25
$ws.Range("K107").Value = 1446.75
$ws.Range("M107").Value = 473.25

$ws.Range("H122").Value = 1728
$ws.Range("I122").Value = 1144
$ws.Range("K122").Value = 3432
$ws.Range("M122").Value = -982

$ws.Range("H132").Value = 1622.4464
$ws.Range("I132").Value = 1068.1904
$ws.Range("J132").Value = 3285.2144
$ws.Range("K132").Value = 3204.5712
$ws.Range("L132").Value = 9855.643199999999
$ws.Range("M132").Value = -674.5711999999999
$ws.Range("N132").Value = -14915.6432
